$d = $word.ActiveDocument

# Position an insertion point right at the very end of the document body
# (just before the final paragraph mark / body end), then splice in the
# two new paragraphs as raw WordprocessingML: a blank paragraph followed
# by a paragraph containing the text "kk".
$end = $d.Content.End - 1
$r = $d.Range($end, $end)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = "<w:p $ns/><w:p $ns><w:r><w:t>kk</w:t></w:r></w:p>"

$r.InsertXML($xml) | Out-Null
